$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty X7/Y7 cells on the existing last row
$ws.Range("X7").Value = 0.22000199999999381
$ws.Range("Y7").Value = "Up"

# Append new row 8 with the new scan's data
$ws.Range("A8").Value = 42649.879120370373
$ws.Range("B8").Value = -2
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 14272
$ws.Range("F8").Value = 855
$ws.Range("G8").Value = 51
$ws.Range("H8").Value = 46
$ws.Range("I8").Value = 85
$ws.Range("J8").Value = 14
$ws.Range("K8").Value = 19350
$ws.Range("L8").Value = 129
$ws.Range("M8").Value = 115
$ws.Range("N8").Value = 53
$ws.Range("O8").Value = 9
$ws.Range("P8").Value = "Named"
$ws.Range("Q8").Value = 47.418521827693588
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = -0.089899999999999994
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("T8").Value = -0.021000000000000001
$ws.Range("T8").NumberFormat = "0.00%"
$ws.Range("U8").Value = 6.65
$ws.Range("V8").Value = 1.88
$ws.Range("W8").Value = -2
